# Insert a new weekly record at the top of the Cereza price series
# (row 46 of the data block) and shift all subsequent rows down by one,
# matching the diff: dimension grows from A1:T102 to A1:T103.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46..102 down to 47..103, leaving a blank row 46 to fill in.
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with the new weekly observation.
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 44589
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100103
$ws.Cells.Item(46, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(46, 9).Value = 100103001
$ws.Cells.Item(46, 10).Value = "Cereza"
$ws.Cells.Item(46, 11).Value = "Santina"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 120
$ws.Cells.Item(46, 14).Value = 5500
$ws.Cells.Item(46, 15).Value = 6000
$ws.Cells.Item(46, 16).Value = 5750
$ws.Cells.Item(46, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(46, 19).Value = 575
$ws.Cells.Item(46, 20).Value = 10
